$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86; existing rows 86.. shift down to 87..
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly price record
$ws.Range("A86").Value = 1
$ws.Range("B86").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C86").Value = "Arica y Parinacota"
$ws.Range("D86").Value = 44711
$ws.Range("E86").Value = 15
$ws.Range("F86").Value = 100112036
$ws.Range("G86").Value = "Caigua"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 130
$ws.Range("K86").Value = 13000
$ws.Range("L86").Value = 14000
$ws.Range("M86").Value = 13500
$ws.Range("N86").Value = "$/caja 20 kilos"
$ws.Range("O86").Value = "Región de Arica y Parinacota"
$ws.Range("P86").Value = 675
$ws.Range("Q86").Value = 20
$ws.Range("R86").Value = "Hortaliza"
